# "Update all city code" — add RPA180..RPA183 Amadeus TJQ office/WSAP rows
# to the Assets sheet (rows 12-19), and drop the now-unused trailing blank
# formatting rows 992:995 at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assets")
$ws.Activate()

# Column A "key" cells for the new WSAP/OfficeID1 rows, entered first
# (mirrors the author's original editing order).
$ws.Range("A13").Value = "RPA180_WSAP"
$ws.Range("A14").Value = "RPA181_OfficeID1"
$ws.Range("A16").Value = "RPA182_OfficeID1"
$ws.Range("A18").Value = "RPA183_OfficeID1"
$ws.Range("A19").Value = "RPA183_WSAP"

# Row 12: Thailand (TH) OfficeID1
$ws.Range("A12").Value = "RPA180_OfficeID1"
$ws.Range("B12").Value = "RPA180_Amadeus_TJQ_TH_OfficeID1"
$ws.Range("D12").Value = "BKKIQ317O"

# Row 14: Vietnam (VN) OfficeID1
$ws.Range("B14").Value = "RPA181_Amadeus_TJQ_VN_OfficeID1"
$ws.Range("D14").Value = "SGNTV2001"

# Row 16: Singapore (SG) OfficeID1
$ws.Range("B16").Value = "RPA182_Amadeus_TJQ_SG_OfficeID1"
$ws.Range("D16").Value = "SINTV2100"

# Row 18: Malaysia (MY) OfficeID1
$ws.Range("B18").Value = "RPA183_Amadeus_TJQ_MY_OfficeID1"
$ws.Range("D18").Value = "KULTV28AA"

# Remaining WSAP "key" cells
$ws.Range("A15").Value = "RPA181_WSAP"
$ws.Range("A17").Value = "RPA182_WSAP"

# WSAP rows reuse the existing "RPA179_Amadeus_TJQ_ID_WSAP" value.
$ws.Range("B13").Value = "RPA179_Amadeus_TJQ_ID_WSAP"
$ws.Range("B15").Value = "RPA179_Amadeus_TJQ_ID_WSAP"
$ws.Range("B17").Value = "RPA179_Amadeus_TJQ_ID_WSAP"
$ws.Range("B19").Value = "RPA179_Amadeus_TJQ_ID_WSAP"

# The sheet had 4 extra trailing blank formatting rows - remove them.
$ws.Rows("992:995").Delete()

# Restore the selection/active cell as left by the author.
$ws.Range("A9").Select()
